# "plus mc Praha 5" - fill in software info for Praha 5 (mestska cast)
# on the overview sheet, and mark the corresponding checkbox columns on
# the "Vybrany software" summary sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Prehledova tabulka
$ws2 = $wb.Worksheets.Item(2)   # Vybrany software

# --- Prehledova tabulka ---------------------------------------------------
# Row 6 = Praha 5: the "ucetnictvi" column used to say "N/A", now it has
# a value.
$ws1.Range("B6").Value = "Ginis – Gordic"

# Row 7 = Praha 6: fill in the newly supplied software details.
$ws1.Range("D7").Value = "Datacentrum2"
$ws1.Range("E7").Value = "Ginis – Gordic"
$ws1.Range("F7").Value = "Materiály RMČ"
$ws1.Range("H7").Value = "VITA"
$ws1.Range("J7").Value = "ISMA"
$ws1.Range("K7").Value = "VERA"
$ws1.Range("L7").Value = "VISA"

# --- Vybrany software (summary) ------------------------------------------
# Row 7 = Praha 5: flag the software columns that now apply.
$ws2.Range("B7").Value = 1
$ws2.Range("D7").Value = 1
$ws2.Range("L7").Value = 1

# --- Selection / active sheet state ---------------------------------------
$ws1.Range("F7").Select()
$ws2.Activate()
$ws2.Range("A8").Select()
